$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C from 45207 to 45208 for all data rows (2-33)
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# Update the hyperlink formulas in row 2 from Logging_SOLLEFTEA to Logging_2283
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/artfynd/A 33036-2023.xlsx", "A 33036-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/kartor/A 33036-2023.png", "A 33036-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/klagomål/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/klagomålsmail/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/tillsyn/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/tillsynsmail/A 33036-2023.docx", "A 33036-2023")'

$wb.Save()
